# Add a "SEAT NO" column in front of the student data.
# This inserts a new column A (shifting the existing STUDENT_NAME /
# MOTHER NAME / SUBn columns one place to the right, B->C, C->D, ...)
# and then fills the new column with the "SEAT NO:" header and the
# per-student seat numbers s1..s8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:G to B:H, inserting a fresh blank column A.
$ws.Columns("A").Insert()

# Header label for the new column (row 4 holds the other column headers).
$ws.Range("A4").Value = "SEAT NO:"

# Seat numbers for each of the 8 students (rows 5-12).
$ws.Range("A5").Value = "s1"
$ws.Range("A6").Value = "s2"
$ws.Range("A7").Value = "s3"
$ws.Range("A8").Value = "s4"
$ws.Range("A9").Value = "s5"
$ws.Range("A10").Value = "s6"
$ws.Range("A11").Value = "s7"
$ws.Range("A12").Value = "s8"

# Match the author's final selection in the saved workbook.
[void]$ws.Range("C9").Select()
